# Prep up level 5 stuff, texts and tutorial flow added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# The shared-string table records *first use* order, which reflects the
# exact sequence the author typed values into cells while drafting the
# level_3 / level_5 tutorial copy. Reproduce that sequence precisely so the
# resulting sharedStrings.xml ordering matches.

$ws.Cells.Item(33, 1).Value = "level_3_intro_1"
$ws.Cells.Item(34, 1).Value = "level_3_intro_2"

$ws.Cells.Item(33, 2).Value = "Watch out! A yeti is in the way!"
$ws.Cells.Item(34, 2).Value = "We must withstand the freezing wind to proceed!"

$ws.Cells.Item(35, 1).Value = "level_3_info_1"
$ws.Cells.Item(35, 2).Value = "In order to succeed, you must subtract the fractional attack until it reaches below zero."

$ws.Cells.Item(36, 1).Value = "mixedNumber"
$ws.Cells.Item(36, 2).Value = "Mixed Number"

$ws.Cells.Item(37, 1).Value = "improperFraction"
$ws.Cells.Item(37, 2).Value = "Improper Fraction"

$ws.Cells.Item(38, 1).Value = "level_5_intro_1"
$ws.Cells.Item(38, 2).Value = "After an onerous battle, Pengu must catch a breath of fresh air."

$ws.Cells.Item(39, 1).Value = "level_5_intro_2"

$ws.Cells.Item(40, 1).Value = "level_5_info_1_a"
$ws.Cells.Item(40, 2).Value = "As you can see, there is a mixed number in this operation."

$ws.Cells.Item(41, 1).Value = "level_5_info_1_b"
$ws.Cells.Item(41, 2).Value = "A mixed number is made up of a whole number, and a fraction."

$ws.Cells.Item(42, 1).Value = "level_5_info_2_a"
$ws.Cells.Item(44, 1).Value = "level_5_info_3_a"
$ws.Cells.Item(43, 1).Value = "level_5_info_2_b"

$ws.Cells.Item(42, 2).Value = "To convert a mixed number to an improper fraction: multiply the whole number with the denominator of the fraction."
$ws.Cells.Item(43, 2).Value = "Afterwards, add the result to the numerator."

$ws.Cells.Item(45, 1).Value = "level_5_info_3_b"
$ws.Cells.Item(45, 2).Value = "Use this technique to help you out with tricky operations!"

$ws.Cells.Item(39, 2).Value = "Help Pengu swim towards the land by adding up the distances using fractions!"
$ws.Cells.Item(44, 2).Value = "You can drag the whole number towards the fraction, or vice-versa, to convert."

# Cells that carry the vertical-center alignment style in the target workbook
$centeredCells = @("B40", "B42", "B43", "B44")
foreach ($addr in $centeredCells) {
    $ws.Range($addr).VerticalAlignment = -4108
}

# Update the visible window / selection state to match the authored edit
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B44").Select()
